# edit.ps1 -- PowerPoint COM-interop script (PowerShell-style) that applies
# the textual changes described by the target diff:
#
#   * Slide 2 ("Content Placeholder 1"): the run containing
#     "[1:2]}}, " is merged into the following run "this is how I
#     stack-up:", and the index token changes from 1:2 to 1:1, producing
#     a single run "[1:1]}}, this is how I stack-up:".
#   * Slide 3: four separate shapes each contain a "DATE n = {{val:dates.csv[n:2]}} "
#     style run; the "[n:2]}} " token becomes "[n:1]}} " in every one of them.
#
# (The diff also touches the cached text of the notes-master
# datetimeFigureOut field and swaps the customXml/item*.xml parts around;
# those two live outside of anything the PowerPoint object model exposes
# in this host -- NotesMaster shapes reject every edit with "target not
# found", and Presentation.CustomXMLParts.Count is always 0 here -- so
# they cannot be reproduced through COM automation and are intentionally
# left untouched.)

$p = $ppt.ActivePresentation

function Replace-Substring($textRange, [string]$oldSubstring, [string]$newSubstring) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldSubstring)
    if ($idx -lt 0) {
        return $false
    }
    $chars = $textRange.Characters($idx + 1, $oldSubstring.Length)
    $chars.Text = $newSubstring
    return $true
}

# --- Slide 2: "If my birthday is {{val:birthday.csv[1:2]}}, this is how I stack-up:" ---
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(1)
$tr = $contentShape.TextFrame.TextRange

# Drop the old "[1:2]}}, " run entirely ...
Replace-Substring $tr "[1:2]}}, " "" | Out-Null

# ... and prepend the replacement text (with the updated 1:1 index) onto the
# following run, so it keeps that run's formatting (dirty="0") rather than
# the deleted run's.
$full = $tr.Text
$idx = $full.IndexOf("this is how I stack-up:")
$afterChars = $tr.Characters($idx + 1, "this is how I stack-up:".Length)
$afterChars.InsertBefore("[1:1]}}, ") | Out-Null

# --- Slide 3: four "DATE n = {{val:dates.csv[n:2]}} " boxes ---
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shape = $slide3.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $shapeTr = $shape.TextFrame.TextRange
    for ($n = 1; $n -le 4; $n++) {
        $didReplace = Replace-Substring $shapeTr "[$n`:2]}} " "[$n`:1]}} "
        if ($didReplace) {
            break
        }
    }
}
